# Add cantrals by cantons
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wipe the existing content/formatting of the sheet so we can rebuild it
# with the new layout (new columns idx/idx2/Name/Date Start/Date End plus
# renamed (MW1)/(MW2)/(GWh) headers).
$ws.Cells.Clear()

# ---- Header row ----
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Style the F1:K1 headers with Arial 9 / general number format (new style)
$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9
$hdr.WrapText = $false

# ---- Data rows ----
# Row 2: Chancy-Pougny
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 509700
$ws.Range("C2").Value = "Chancy-Pougny"
$ws.Range("D2").Value = 1925
$ws.Range("F2").Value = 520
$ws.Range("G2").Value = 26.16
$ws.Range("H2").Value = 24.87
$ws.Range("I2").Value = 58.14
$ws.Range("J2").Value = 77.52
$ws.Range("K2").Value = 135.66

# Row 3: Verbois
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 509600
$ws.Range("C3").Value = "Verbois"
$ws.Range("D3").Value = 1943
$ws.Range("E3").Value = 1999
$ws.Range("F3").Value = 620
$ws.Range("G3").Value = 102.8
$ws.Range("H3").Value = 98
$ws.Range("I3").Value = 211
$ws.Range("J3").Value = 255
$ws.Range("K3").Value = 466

# Row 4: Seujet
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 509450
$ws.Range("C4").Value = "Seujet"
$ws.Range("D4").Value = 1994
$ws.Range("F4").Value = 405
$ws.Range("G4").Value = 8.6999999999999993
$ws.Range("H4").Value = 5.6
$ws.Range("I4").Value = 9.8000000000000007
$ws.Range("J4").Value = 10.199999999999999
$ws.Range("K4").Value = 20

# Apply fonts / number formats matching the workbook's existing styles:
#  - A,B (idx/idx2) -> Arial 9, integer "0" format
#  - C (Name) -> Arial 9, general format
#  - D (Date Start) -> Arial 9, integer "0" format
#  - E3 (Date End, only Verbois has one) -> Arial 9, integer "0" format
#  - F:K (data) -> Arial 9, "0.00" format
# NOTE: number-format/font is only ever applied to ranges whose every cell
# already has a value - applying it to a range that also covers blank cells
# would materialize phantom empty <c> elements that aren't present in the
# target workbook.
$colAB = $ws.Range("A2:B4")
$colAB.Font.Name = "Arial"
$colAB.Font.Size = 9
$colAB.NumberFormat = "0"

$colC = $ws.Range("C2:C4")
$colC.Font.Name = "Arial"
$colC.Font.Size = 9

$colD = $ws.Range("D2:D4")
$colD.Font.Name = "Arial"
$colD.Font.Size = 9
$colD.NumberFormat = "0"

$colE = $ws.Range("E3")
$colE.Font.Name = "Arial"
$colE.Font.Size = 9
$colE.NumberFormat = "0"

$dataCols = $ws.Range("F2:K4")
$dataCols.Font.Name = "Arial"
$dataCols.Font.Size = 9
$dataCols.NumberFormat = "0.00"

# ---- Selection state ----
[void]$ws.Range("I18").Select()
